$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data appended below the existing rows (row 16)
$rowNum = 16

# Column A is a text timestamp, not a date value, matching the other rows.
$ws.Cells.Item($rowNum, 1).Value = "2024-09-04 16:28:43"

$values = @(0, 0, 0, 0, 25, 100, 0, 0, 100, 0, 0, 50, 0, 33.33333333333333, 20, 100, 0, 0, 100, 0, 0, 0, 100, 0, 0, 100)

$col = 2
foreach ($v in $values) {
    $ws.Cells.Item($rowNum, $col).Value = $v
    $col = $col + 1
}
